$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# 1) Replace the ID placeholder text in the first paragraph's first run.
$p1.Range.Find.Execute("**ID__AFFARS_5325_topic_22__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5325_7703_2__ID**", 2)

# 2) Remove the trailing run that holds just a single space " " in paragraph 1.
$p1.Range.Find.Execute(" ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3) Update the left indent on paragraph 1 from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# 4) Add a paragraph border (space = 5pt on each side) to paragraph 1.
$borders = $p1.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
